# Clean up the "Reference" column (A) by stripping the stray trailing "16"
# digits that had been accidentally appended to each Bible verse reference
# (e.g. "Malachi 1:116" -> "Malachi 1:1"), making the references human
# readable again. Column B (Text) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$references = @{
    2  = "Malachi 1:1"
    3  = "Malachi 1:3"
    4  = "Malachi 1:4"
    5  = "Malachi 1:9"
    6  = "Malachi 1:11"
    7  = "Malachi 1:12"
    8  = "Malachi 1:13"
    9  = "Malachi 1:14"
    10 = "Malachi 2:1"
    11 = "Malachi 2:3"
    12 = "Malachi 2:4"
    13 = "Malachi 2:5"
    14 = "Malachi 2:7"
    15 = "Malachi 2:8"
    16 = "Malachi 2:9"
    17 = "Malachi 2:10"
    18 = "Malachi 2:11"
    19 = "Malachi 2:12"
    20 = "Malachi 2:13"
    21 = "Malachi 2:16"
    22 = "Malachi 3:2"
    23 = "Malachi 3:3"
    24 = "Malachi 3:5"
    25 = "Malachi 3:6"
    26 = "Malachi 3:8"
    27 = "Malachi 3:11"
    28 = "Malachi 3:12"
    29 = "Malachi 3:14"
    30 = "Malachi 3:15"
    31 = "Malachi 3:16"
    32 = "Malachi 3:17"
    33 = "Malachi 3:18"
    34 = "Malachi 4:1"
    35 = "Malachi 4:2"
    36 = "Malachi 4:3"
    37 = "Malachi 4:6"
}

foreach ($row in $references.Keys) {
    $ws.Cells.Item($row, 1).Value = $references[$row]
}
